$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "955.993.998.1001.1006.1009.10424.20103.20125.20310.20384.40151.50623.60159.60162.60192.60225.70103.70113.10425.60126.40139.20205.60280"

$ws.Range("B10").Select()
